$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values. Cells whose new text looks like a number (e.g. "1.005")
# are first formatted as Text so Excel keeps the exact original string instead of
# auto-converting it into a numeric value.

$ws.Range('D2').Value = '27.952.29'
$ws.Range('D3').Value = '1.903.26'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.05'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4819'
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3794'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07367'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9312'
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.75'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07743'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '1.941.03'
$ws.Range('E13').Value = '  +4.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.479'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.624'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.66'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008868'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').Value = '27.987.58'
$ws.Range('E20').Value = '  +1.86%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.67'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.149'
$ws.Range('D23').Value = '2.159.55'
$ws.Range('E23').Value = '  +3.05%  '
$ws.Range('E24').Value = '  +2.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '156.01'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.907'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.46'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.125'
$ws.Range('E28').Value = '  +5.98%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '117.03'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.962'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08937'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.262'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.251'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7668'
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.665'
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02056'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.553'
$ws.Range('E37').Value = '  -5.69%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.105'
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5485'
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05275'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.995'
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.939'
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1526'
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.476'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '109.89'
$ws.Range('E45').Value = '  +6.54%  '
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4805'
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.643'
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '67.80'
$ws.Range('E50').Value = '  +0.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06079'
$ws.Range('E51').Value = '  -0.41%  '
